$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 102 and 103 (columns F:V only; A:E stay as-is) ---
$row102 = @("Estudiantes L.P.", 0, "Godoy Cruz", 1, 1.86, "02/10/2023 22:42", 2.17, "07/10/2023 23:59", 3.37, "02/10/2023 22:42", 3.06, "07/10/2023 23:59", 4.33, "02/10/2023 22:42", 4.04, "07/10/2023 23:59", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/estudiantes-l-p-godoy-cruz/x6gob2XD/")
$row103 = @("Lanus", 0, "Defensa y Justicia", 2, 1.81, "01/10/2023 02:42", 1.9, "07/10/2023 23:58", 3.58, "01/10/2023 02:42", 3.47, "07/10/2023 23:58", 4.7, "01/10/2023 02:42", 4.55, "07/10/2023 23:58", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/lanus-defensa-y-justicia/EX7jxLAf/")

for ($i = 0; $i -lt $row102.Length; $i++) {
    $col = 6 + $i
    $ws.Cells.Item(102, $col).Value = $row102[$i]
    $ws.Cells.Item(103, $col).Value = $row103[$i]
}

# --- Swap rows 110 and 111 (columns F:V only; A:E stay as-is) ---
$row110 = @("Arsenal Sarandi", 0, "Banfield", 0, 3.26, "02/10/2023 22:12", 3.55, "09/10/2023 22:59", 2.96, "02/10/2023 22:12", 2.9, "09/10/2023 22:58", 2.51, "02/10/2023 22:12", 2.46, "09/10/2023 22:53", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/arsenal-sarandi-banfield/lUOSMuB7/")
$row111 = @("Barracas Central", 2, "Colon Santa Fe", 1, 2.5, "02/10/2023 22:42", 2.6, "09/10/2023 22:36", 3.13, "02/10/2023 22:42", 2.84, "09/10/2023 22:36", 2.93, "02/10/2023 22:42", 2.91, "09/10/2023 22:36", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/barracas-central-colon-santa-fe/4x7Fpq3E/")

for ($i = 0; $i -lt $row110.Length; $i++) {
    $col = 6 + $i
    $ws.Cells.Item(110, $col).Value = $row110[$i]
    $ws.Cells.Item(111, $col).Value = $row111[$i]
}

# --- Append new row 183 (copy formatting/style from row 182, then set values) ---
$ws.Range("A182:V182").Copy($ws.Range("A183:V183"))

$newRow = @(182, "argentina", "copa-de-la-liga-profesional", "2023", 45244.04166666666, "Instituto", 0, "Barracas Central", 0, 1.71, "07/11/2023 05:11", 1.73, "14/11/2023 00:34", 3.58, "07/11/2023 05:11", 3.38, "14/11/2023 00:33", 5.77, "07/11/2023 05:11", 6.16, "14/11/2023 00:48", "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/instituto-barracas-central/vXbZJXeJ/")

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $col = 1 + $i
    $ws.Cells.Item(183, $col).Value = $newRow[$i]
}

# Column D ("temporada") is stored as text ("2023"), not a number, in every
# other row. A plain string assignment of a purely-numeric-looking string
# gets auto-coerced to a number, so force text via NumberFormat and then
# restore the original (default) cell style copied from row 182's D column.
$ws.Range("D183").NumberFormat = "@"
$ws.Range("D183").Value = "2023"
$ws.Range("D183").Style = $ws.Range("D182").Style
